$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.454345666666666
$ws.Range("H2").Value = 7.363036999999999
$ws.Range("I2").Value = 0.2857469401257222
$ws.Range("J2").Value = 0.3027613480760606
$ws.Range("M2").Value = 12.384602
$ws.Range("N2").Value = 37.153806
$ws.Range("O2").Value = 0.9505912801254437
$ws.Range("P2").Value = 0.9566996763338882
$ws.Range("Q2").Value = 30.39609425209133
$ws.Range("R2").Value = 273.564848268822
$ws.Range("S2").Value = 0.2716285496060387
$ws.Range("T2").Value = 0.2896516837107788
$ws.Range("G3").Value = 2.454345666666666
$ws.Range("H3").Value = 7.363036999999999
$ws.Range("I3").Value = 0.2857469401257222
$ws.Range("J3").Value = 0.3027613480760606
$ws.Range("O3").Value = 0.01115977065643923
$ws.Range("P3").Value = 0.01123148212927739
$ws.Range("Q3").Value = 0.3568446795136666
$ws.Range("R3").Value = 3.211602115623
$ws.Range("S3").Value = 0.003188870317582332
$ws.Range("T3").Value = 0.003400458670352205
$ws.Range("G4").Value = 2.454345666666666
$ws.Range("H4").Value = 7.363036999999999
$ws.Range("I4").Value = 0.2857469401257222
$ws.Range("J4").Value = 0.3027613480760606
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.153231
$ws.Range("N4").Value = 0.459693
$ws.Range("O4").Value = 0.01176138340536917
$ws.Range("P4").Value = 0.01183696077631869
$ws.Range("Q4").Value = 0.3760818408489999
$ws.Range("R4").Value = 3.384736567640999
$ws.Range("S4").Value = 0.003360779319729687
$ws.Range("T4").Value = 0.003583774201761699
$ws.Range("G5").Value = 2.454345666666666
$ws.Range("H5").Value = 7.363036999999999
$ws.Range("I5").Value = 0.2857469401257222
$ws.Range("J5").Value = 0.3027613480760606
$ws.Range("M5").Value = 0.249552
$ws.Range("N5").Value = 0.499104
$ws.Range("O5").Value = 0.01915458850739529
$ws.Range("P5").Value = 0.01285178254031226
$ws.Range("Q5").Value = 0.6124868698079999
$ws.Range("R5").Value = 3.674921218848
$ws.Range("S5").Value = 0.005473365055355527
$ws.Range("T5").Value = 0.003891023007085318
$ws.Range("G6").Value = 2.454345666666666
$ws.Range("H6").Value = 7.363036999999999
$ws.Range("I6").Value = 0.2857469401257222
$ws.Range("J6").Value = 0.3027613480760606
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09553633333333333
$ws.Range("N6").Value = 0.286609
$ws.Range("O6").Value = 0.0073329773053526
$ws.Range("P6").Value = 0.007380098220203317
$ws.Range("Q6").Value = 0.2344791857258889
$ws.Range("R6").Value = 2.110312671533
$ws.Range("S6").Value = 0.002095375827015869
$ws.Range("T6").Value = 0.002234408486082492
$ws.Range("G7").Value = 4.686805000000001
$ws.Range("I7").Value = 0.5456607868665887
$ws.Range("J7").Value = 0.5781514068052169
$ws.Range("M7").Value = 12.384602
$ws.Range("N7").Value = 37.153806
$ws.Range("O7").Value = 0.9505912801254437
$ws.Range("P7").Value = 0.9566996763338882
$ws.Range("Q7").Value = 58.04421457661001
$ws.Range("R7").Value = 522.3979311894901
$ws.Range("S7").Value = 0.5187003859017675
$ws.Range("T7").Value = 0.5531172637625331
$ws.Range("G8").Value = 4.686805000000001
$ws.Range("I8").Value = 0.5456607868665887
$ws.Range("J8").Value = 0.5781514068052169
$ws.Range("O8").Value = 0.01115977065643923
$ws.Range("P8").Value = 0.01123148212927739
$ws.Range("Q8").Value = 0.6814286393650001
$ws.Range("R8").Value = 6.132857754285
$ws.Range("S8").Value = 0.006089449237643299
$ws.Range("T8").Value = 0.006493497193549374
$ws.Range("G9").Value = 4.686805000000001
$ws.Range("I9").Value = 0.5456607868665887
$ws.Range("J9").Value = 0.5781514068052169
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.153231
$ws.Range("N9").Value = 0.459693
$ws.Range("O9").Value = 0.01176138340536917
$ws.Range("P9").Value = 0.01183696077631869
$ws.Range("Q9").Value = 0.718163816955
$ws.Range("R9").Value = 6.463474352595
$ws.Range("S9").Value = 0.00641772572361338
$ws.Range("T9").Value = 0.006843555525126823
$ws.Range("G10").Value = 4.686805000000001
$ws.Range("I10").Value = 0.5456607868665887
$ws.Range("J10").Value = 0.5781514068052169
$ws.Range("M10").Value = 0.249552
$ws.Range("N10").Value = 0.499104
$ws.Range("O10").Value = 0.01915458850739529
$ws.Range("P10").Value = 0.01285178254031226
$ws.Range("Q10").Value = 1.16960156136
$ws.Range("R10").Value = 7.01760936816
$ws.Range("S10").Value = 0.01045190783705103
$ws.Range("T10").Value = 0.007430276155636257
$ws.Range("G11").Value = 4.686805000000001
$ws.Range("I11").Value = 0.5456607868665887
$ws.Range("J11").Value = 0.5781514068052169
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.09553633333333333
$ws.Range("N11").Value = 0.286609
$ws.Range("O11").Value = 0.0073329773053526
$ws.Range("P11").Value = 0.007380098220203317
$ws.Range("Q11").Value = 0.4477601647483334
$ws.Range("R11").Value = 4.029841482735001
$ws.Range("S11").Value = 0.004001318166513537
$ws.Range("T11").Value = 0.004266814168371226
$ws.Range("G12").Value = 1.4480775
$ws.Range("H12").Value = 2.896155
$ws.Range("I12").Value = 0.1685922730076891
$ws.Range("J12").Value = 0.1190872451187225
$ws.Range("M12").Value = 12.384602
$ws.Range("N12").Value = 37.153806
$ws.Range("O12").Value = 0.9505912801254437
$ws.Range("P12").Value = 0.9566996763338882
$ws.Range("Q12").Value = 17.933863502655
$ws.Range("R12").Value = 107.60318101593
$ws.Range("S12").Value = 0.1602623446176375
$ws.Range("T12").Value = 0.1139307288605762
$ws.Range("G13").Value = 1.4480775
$ws.Range("H13").Value = 2.896155
$ws.Range("I13").Value = 0.1685922730076891
$ws.Range("J13").Value = 0.1190872451187225
$ws.Range("O13").Value = 0.01115977065643923
$ws.Range("P13").Value = 0.01123148212927739
$ws.Range("Q13").Value = 0.2105403319575
$ws.Range("R13").Value = 1.263241991745
$ws.Range("S13").Value = 0.001881451101213601
$ws.Range("T13").Value = 0.001337526265375807
$ws.Range("G14").Value = 1.4480775
$ws.Range("H14").Value = 2.896155
$ws.Range("I14").Value = 0.1685922730076891
$ws.Range("J14").Value = 0.1190872451187225
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.153231
$ws.Range("N14").Value = 0.459693
$ws.Range("O14").Value = 0.01176138340536917
$ws.Range("P14").Value = 0.01183696077631869
$ws.Range("Q14").Value = 0.2218903634025
$ws.Range("R14").Value = 1.331342180415
$ws.Range("S14").Value = 0.001982878362026104
$ws.Range("T14").Value = 0.001409631049430168
$ws.Range("G15").Value = 1.4480775
$ws.Range("H15").Value = 2.896155
$ws.Range("I15").Value = 0.1685922730076891
$ws.Range("J15").Value = 0.1190872451187225
$ws.Range("M15").Value = 0.249552
$ws.Range("N15").Value = 0.499104
$ws.Range("O15").Value = 0.01915458850739529
$ws.Range("P15").Value = 0.01285178254031226
$ws.Range("Q15").Value = 0.36137063628
$ws.Range("R15").Value = 1.44548254512
$ws.Range("S15").Value = 0.003229315614988731
$ws.Range("T15").Value = 0.001530483377590684
$ws.Range("G16").Value = 1.4480775
$ws.Range("H16").Value = 2.896155
$ws.Range("I16").Value = 0.1685922730076891
$ws.Range("J16").Value = 0.1190872451187225
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.09553633333333333
$ws.Range("N16").Value = 0.286609
$ws.Range("O16").Value = 0.0073329773053526
$ws.Range("P16").Value = 0.007380098220203317
$ws.Range("Q16").Value = 0.1383440147325
$ws.Range("R16").Value = 0.8300640883949999
$ws.Range("S16").Value = 0.001236283311823194
$ws.Range("T16").Value = 0.0008788755657496
